# "upgrade to Revit 2016" -- bump the stale 2014 date stamps to 2015.
$p = $ppt.ActivePresentation

# 1) Handout master: the auto date placeholder was last cached as
#    4/22/2014 -- refresh it to 2/3/2015.
$handoutDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$handoutDate.Text = "2/3/2015"

# 2) Notes master: same stale auto date placeholder.
$notesDate = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDate.Text = "2/3/2015"

# 3) Slide 2 ("About these materials ...") has a right-aligned
#    "March 2014 " line near the bottom of the body placeholder;
#    bump the year run from 2014 to 2015 without touching the rest.
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shape = $slide2.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $shapeRange = $shape.TextFrame.TextRange
        $shapeText = $shapeRange.Text
        $yearIndex = $shapeText.IndexOf("2014 ")
        if ($yearIndex -ge 0) {
            $yearRange = $shapeRange.Characters($yearIndex + 1, 5)
            $yearRange.Text = "2015 "
        }
    }
}
